$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 201. This shifts the existing rows 201-234
# down to 202-235, preserving all of their data and formatting.
$ws.Rows.Item(201).Insert()

# Populate the new row 201 with the new record's data.
$ws.Cells.Item(201, 1).Value = 4
$ws.Cells.Item(201, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(201, 3).Value = "Los Lagos"
$ws.Cells.Item(201, 4).Value = 44694
$ws.Cells.Item(201, 5).Value = 10
$ws.Cells.Item(201, 6).Value = 100112032
$ws.Cells.Item(201, 7).Value = "Zapallo italiano"
$ws.Cells.Item(201, 8).Value = "Sin especificar"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 200
$ws.Cells.Item(201, 11).Value = 22000
$ws.Cells.Item(201, 12).Value = 22000
$ws.Cells.Item(201, 13).Value = 22000
$ws.Cells.Item(201, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(201, 15).Value = "Región Metropolitana"
$ws.Cells.Item(201, 16).Value = 440
$ws.Cells.Item(201, 17).Value = 50
$ws.Cells.Item(201, 18).Value = "Hortaliza"
